$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 0.794582
$ws.Range("N2").Value = 2.383746
$ws.Range("O2").Value = 0.03449752952410986
$ws.Range("P2").Value = 0.03449752952410985
$ws.Range("Q2").Value = 0.04373246897666667
$ws.Range("R2").Value = 0.39359222079
$ws.Range("S2").Value = 0.03449752952410986
$ws.Range("T2").Value = 0.03449752952410985

# Row 3 (only O, P, S, T change)
$ws.Range("O3").Value = 0.8945489325574519
$ws.Range("P3").Value = 0.8945489325574517
$ws.Range("S3").Value = 0.8945489325574519
$ws.Range("T3").Value = 0.8945489325574517

# Row 4
$ws.Range("M4").Value = 0.2871986666666667
$ws.Range("N4").Value = 0.8615959999999999
$ws.Range("O4").Value = 0.01246900191876775
$ws.Range("P4").Value = 0.01246900191876775
$ws.Range("Q4").Value = 0.01580693594888889
$ws.Range("R4").Value = 0.14226242354
$ws.Range("S4").Value = 0.01246900191876775
$ws.Range("T4").Value = 0.01246900191876775

# Row 5
$ws.Range("M5").Value = 1.149534666666667
$ws.Range("N5").Value = 3.448604
$ws.Range("O5").Value = 0.04990813547540859
$ws.Range("P5").Value = 0.04990813547540859
$ws.Range("Q5").Value = 0.06326847216222221
$ws.Range("R5").Value = 0.56941624946
$ws.Range("S5").Value = 0.04990813547540859
$ws.Range("T5").Value = 0.04990813547540859

# Row 6
$ws.Range("M6").Value = 0.1975403333333333
$ws.Range("N6").Value = 0.5926210000000001
$ws.Range("O6").Value = 0.008576400524262026
$ws.Range("P6").Value = 0.008576400524262026
$ws.Range("Q6").Value = 0.01087229071277778
$ws.Range("R6").Value = 0.09785061641500002
$ws.Range("S6").Value = 0.008576400524262026
$ws.Range("T6").Value = 0.008576400524262026
